# DG: update model API
#
# The sequence-diagram textbox that shows the AddressBook model call
# "deletePerson(p)" is being renamed to "deletePersons(p)" (the model
# API now accepts/returns a list). The call is split across two runs
# in the original deck:
#   run1 = "deletePerson"   (formatting A: dirty/err/smtClean)
#   run2 = "(p)"            (formatting B: dirty/smtClean)
# and after the edit it must read:
#   run1 = "deletePersons(p"
#   run2 = ")"
# i.e. only an "s" is inserted after "Person" and the ")" moves to its
# own trailing run - each run keeps its own original character
# formatting (color etc.), so we edit each run's characters in place
# instead of overwriting the whole TextRange (which would merge the
# runs and lose the per-run formatting split).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $null
try {
    $shp = $s.Shapes.Item("TextBox 28")
} catch {
    $shp = $null
}

if (($shp -eq $null) -or ($shp.TextFrame.TextRange.Text -ne "deletePerson(p)")) {
    # Fall back to locating the shape by its current text, in case the
    # shape name ever differs from "TextBox 28".
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.HasTextFrame -and ($candidate.TextFrame.TextRange.Text -eq "deletePerson(p)")) {
            $shp = $candidate
        }
    }
}

if ($shp -ne $null) {
    $tr = $shp.TextFrame.TextRange
    if ($tr.Text -eq "deletePerson(p)") {
        # Edit the trailing run ("(p)" -> ")") first while its character
        # offsets (13..15) are still valid, then grow the leading run
        # ("deletePerson" -> "deletePersons(p"). Doing it in this order
        # means the first edit never shifts indices out from under the
        # second one.
        $tr.Characters(13, 3).Text = ")"
        $tr.Characters(1, 12).Text = "deletePersons(p"
    }
}
